$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.284000000000002
$ws.Range("C3").Value = -12.721
$ws.Range("E3").Value = 16.736
$ws.Range("E12").Value = 17.877
$ws.Range("C14").Value = -12.509
$ws.Range("C21").Value = -12.499
$ws.Range("C23").Value = -12.594
$ws.Range("E24").Value = 17.108
$ws.Range("C25").Value = -11.582
$ws.Range("D25").Value = -7.481
$ws.Range("E25").Value = 16.901
$ws.Range("C26").Value = -13.139
$ws.Range("D27").Value = -8.467000000000002
$ws.Range("C29").Value = -12.2
$ws.Range("D31").Value = -8.301
$ws.Range("D39").Value = -7.846000000000001
$ws.Range("D48").Value = -7.423
$ws.Range("E50").Value = 16.583
$ws.Range("D51").Value = -8.373999999999999
$ws.Range("D52").Value = -7.423
$ws.Range("C53").Value = -12.193
$ws.Range("E53").Value = 17.161
$ws.Range("D55").Value = -8.065000000000001
$ws.Range("D56").Value = -8.122
$ws.Range("C57").Value = -13.401
$ws.Range("D57").Value = -8.537000000000001
$ws.Range("E57").Value = 16.687
$ws.Range("C59").Value = -13.063
$ws.Range("E61").Value = 16.628
$ws.Range("E63").Value = 17.609
$ws.Range("C69").Value = -10.683
$ws.Range("E70").Value = 17.547
$ws.Range("D73").Value = -8.004000000000001
$ws.Range("C79").Value = -12.491
$ws.Range("C83").Value = -13.169
$ws.Range("E86").Value = 16.655
$ws.Range("D89").Value = -6.884
$ws.Range("D90").Value = -7.442
$ws.Range("C91").Value = -10.621
$ws.Range("D92").Value = -6.863
$ws.Range("C93").Value = -11.886
$ws.Range("E98").Value = 16.494
$ws.Range("E100").Value = 16.805
$ws.Range("E102").Value = 16.49

"Applied 44 cell updates"
